# Realestate Update resale numbers 2024-01-28 11:28
# Appends a new data row (row 99) to the CityResaleNum sheet with the
# 2024-01-28 11:28:50 resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 99

# --- Text columns (Date / Time / Weekday / Week) -----------------------
# Assigning a plain string to .Value lets Excel auto-infer a date/number
# type (e.g. "2024-01-28" -> a date serial, "04" -> 4), which would lose
# the original literal text. Forcing NumberFormat to "@" (Text) before
# the assignment keeps the literal string, and ClearFormats() afterwards
# drops back to the default (unstyled) cell format so no stray explicit
# "General" style gets attached to the cell.
$c = $ws.Cells.Item($row, 1)
$c.NumberFormat = "@"
$c.Value = "2024-01-28"
$c.ClearFormats()

$ws.Cells.Item($row, 2).Value = "11:28:50"
$ws.Cells.Item($row, 3).Value = "Sunday"

$c = $ws.Cells.Item($row, 4)
$c.NumberFormat = "@"
$c.Value = "04"
$c.ClearFormats()

# --- Numeric columns (city resale numbers) ------------------------------
$ws.Cells.Item($row, 5).Value = 137323
$ws.Cells.Item($row, 6).Value = 141829
$ws.Cells.Item($row, 7).Value = 171212
$ws.Cells.Item($row, 8).Value = 149115
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 121427
$ws.Cells.Item($row, 11).Value = 223856
$ws.Cells.Item($row, 12).Value = 257213
$ws.Cells.Item($row, 13).Value = 185450
$ws.Cells.Item($row, 14).Value = 110034
$ws.Cells.Item($row, 15).Value = 41405
$ws.Cells.Item($row, 16).Value = 30828
$ws.Cells.Item($row, 17).Value = 73627
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42121
$ws.Cells.Item($row, 20).Value = -1
